$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to be treated as text so that
# numeric-looking price/percentage strings are preserved exactly as text,
# matching the original inline-string cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value2 = '29.220.54'
$ws.Range("E2").Value2 = '  -0.06%  '
$ws.Range("D3").Value2 = '1.826.69'
$ws.Range("E3").Value2 = '  +0.03%  '
$ws.Range("D4").Value2 = '1.003'
$ws.Range("E4").Value2 = '  +0.15%  '
$ws.Range("D5").Value2 = '235.17'
$ws.Range("E5").Value2 = '  -0.53%  '
$ws.Range("D6").Value2 = '0.6001'
$ws.Range("E6").Value2 = '  -0.59%  '
$ws.Range("D7").Value2 = '1.005'
$ws.Range("E7").Value2 = '  +0.19%  '
$ws.Range("D8").Value2 = '0.06949'
$ws.Range("E8").Value2 = '  -2.89%  '
$ws.Range("D9").Value2 = '0.2764'
$ws.Range("E9").Value2 = '  -1.86%  '
$ws.Range("D10").Value2 = '23.46'
$ws.Range("E10").Value2 = '  -3.09%  '
$ws.Range("D11").Value2 = '0.07609'
$ws.Range("E11").Value2 = '  -0.85%  '
$ws.Range("D12").Value2 = '1.833.71'
$ws.Range("E12").Value2 = '  +0.91%  '
$ws.Range("D13").Value2 = '4.735'
$ws.Range("E13").Value2 = '  -1.26%  '
$ws.Range("D14").Value2 = '0.6300'
$ws.Range("E14").Value2 = '  -2.23%  '
$ws.Range("D15").Value2 = '0.000009834'
$ws.Range("E15").Value2 = '  +0.36%  '
$ws.Range("D16").Value2 = '77.48'
$ws.Range("E16").Value2 = '  -2.84%  '
$ws.Range("D17").Value2 = '29.019.36'
$ws.Range("E17").Value2 = '  -0.61%  '
$ws.Range("D18").Value2 = '5.549'
$ws.Range("E18").Value2 = '  -8.46%  '
$ws.Range("D19").Value2 = '216.09'
$ws.Range("E19").Value2 = '  -6.38%  '
$ws.Range("D20").Value2 = '1.004'
$ws.Range("E20").Value2 = '  +0.17%  '
$ws.Range("D21").Value2 = '11.56'
$ws.Range("E21").Value2 = '  -2.18%  '
$ws.Range("D22").Value2 = '6.863'
$ws.Range("E22").Value2 = '  -2.52%  '
$ws.Range("D23").Value2 = '1.005'
$ws.Range("E23").Value2 = '  +0.18%  '
$ws.Range("D24").Value2 = '155.96'
$ws.Range("E24").Value2 = '  -0.17%  '
$ws.Range("D25").Value2 = '7.948'
$ws.Range("E25").Value2 = '  -2.24%  '
$ws.Range("E26").Value2 = '  +0.25%  '
$ws.Range("E27").Value2 = '  -1.21%  '
$ws.Range("D28").Value2 = '0.06426'
$ws.Range("E28").Value2 = '  -5.65%  '
$ws.Range("D29").Value2 = '1.416'
$ws.Range("E29").Value2 = '  -3.38%  '
$ws.Range("E30").Value2 = '  -1.30%  '
$ws.Range("D31").Value2 = '3.823'
$ws.Range("E31").Value2 = '  +1.24%  '
$ws.Range("D32").Value2 = '3.782'
$ws.Range("E32").Value2 = '  -1.87%  '
$ws.Range("D33").Value2 = '1.094'
$ws.Range("E33").Value2 = '  -3.71%  '
$ws.Range("D34").Value2 = '1.723'
$ws.Range("E34").Value2 = '  +0.07%  '
$ws.Range("D35").Value2 = '0.6465'
$ws.Range("E35").Value2 = '  -2.52%  '
$ws.Range("D36").Value2 = '2.544'
$ws.Range("E36").Value2 = '  +0.58%  '
$ws.Range("D37").Value2 = '2.755'
$ws.Range("E37").Value2 = '  -0.25%  '
$ws.Range("D38").Value2 = '0.01756'
$ws.Range("E38").Value2 = '  -0.82%  '
$ws.Range("D39").Value2 = '6.598'
$ws.Range("E39").Value2 = '  +0.16%  '
$ws.Range("D40").Value2 = '1.136.61'
$ws.Range("E40").Value2 = '  -7.21%  '
$ws.Range("D41").Value2 = '0.8933'
$ws.Range("E41").Value2 = '  -3.53%  '
$ws.Range("E42").Value2 = '  +0.14%  '
$ws.Range("D43").Value2 = '1.996.99'
$ws.Range("E43").Value2 = '  +1.80%  '
$ws.Range("D44").Value2 = '100.82'
$ws.Range("E44").Value2 = '  +1.39%  '
$ws.Range("E45").Value2 = '  -2.22%  '
$ws.Range("D46").Value2 = '0.00000000114'
$ws.Range("E46").Value2 = '  -3.65%  '
$ws.Range("D47").Value2 = '1.618'
$ws.Range("E47").Value2 = '  -1.24%  '
$ws.Range("D48").Value2 = '8.481'
$ws.Range("E48").Value2 = '  -0.62%  '
$ws.Range("D49").Value2 = '0.05504'
$ws.Range("E49").Value2 = '  -1.71%  '
$ws.Range("D50").Value2 = '0.4534'
$ws.Range("E50").Value2 = '  -0.75%  '
$ws.Range("D51").Value2 = '6.378'
$ws.Range("E51").Value2 = '  -3.62%  '
